$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on column D (Price) cells so numeric-looking strings
# are retained as text, matching the original inlineStr cell type, then reset
# the style back to Normal so no stray formatting differences are introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.748.22"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.292.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +18.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.51%  "

$ws.Range("E13").Value = "  +0.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.635.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.845"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.292.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.642.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000109"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.65%  "

$ws.Range("E20").Value = "  +4.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.28%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.24%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.86%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.36%  "

$ws.Range("E29").Value = "  -1.91%  "

$ws.Range("E30").Value = "  -1.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "176.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0929"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.59%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.36%  "

$ws.Range("E34").Value = "  +3.53%  "

$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.34%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0361"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.109"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.243"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.48%  "

$ws.Range("E45").Value = "  +5.50%  "

$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.29%  "

$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("E50").Value = "  +3.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.450"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.94%  "
